$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Metadata sheet updates
$wsMeta.Range("B4").Value = "QCEthnicityCodeSystem"
$wsMeta.Range("B8").Value = "2025-09-23T20:31:36+00:00"

# Concepts sheet updates (Display column B, Code column C share the same text)
$wsConcepts.Range("B4").Value = "First Nation, Inuit, Metis"
$wsConcepts.Range("C4").Value = "First Nation, Inuit, Metis"

$wsConcepts.Range("B5").Value = "European"
$wsConcepts.Range("C5").Value = "European"

$wsConcepts.Range("B6").Value = "Arab"
$wsConcepts.Range("C6").Value = "Arab"

$wsConcepts.Range("B7").Value = "Latin American"
$wsConcepts.Range("C7").Value = "Latin American"
